# Updates with response to Jason's comments:
#  - Column L ("color") on "slide.pot.objects" gets a formula-based value
#    of ="199,183,199" for the rows that should stay visible, and an
#    AutoFilter is applied on column L (index 12) restricted to that
#    value, which hides every other data row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(3)

# --- Set the "keep" rows' L column to the new shared formula value ---
# (grouped into contiguous ranges so Excel creates shared formulas the
# same way the original author's edit did)
$ws.Range("L2:L3").Formula = '="199,183,199"'
$ws.Range("L2:L3").NumberFormat = "#,##0"

$ws.Range("L5:L6").Formula = '="199,183,199"'
$ws.Range("L5:L6").NumberFormat = "#,##0"

$ws.Range("L8").Formula = '="199,183,199"'
$ws.Range("L8").NumberFormat = "#,##0"

$ws.Range("L12").Formula = '="199,183,199"'
$ws.Range("L12").NumberFormat = "#,##0"

$ws.Range("L19:L20").Formula = '="199,183,199"'
$ws.Range("L19:L20").NumberFormat = "#,##0"

$ws.Range("L23").Formula = '="199,183,199"'
$ws.Range("L23").NumberFormat = "#,##0"

$ws.Range("L27:L28").Formula = '="199,183,199"'
$ws.Range("L27:L28").NumberFormat = "#,##0"

$ws.Range("L30:L31").Formula = '="199,183,199"'
$ws.Range("L30:L31").NumberFormat = "#,##0"

$ws.Range("L38:L39").Formula = '="199,183,199"'
$ws.Range("L38:L39").NumberFormat = "#,##0"

# --- Apply the AutoFilter on column L (12th column) to only the
#     "199,183,199" value; this hides the rest of the data rows ---
$ws.Range("A1:O40").AutoFilter(12, @("199,183,199"), 7)

# --- Update the frozen-pane scroll position / active selection ---
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 2
$ws.Range("L6").Select()

$wb.Save()
